$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44526
$ws.Range("J2").Value = 100
$ws.Range("D3").Value = 44489
$ws.Range("K3").Value = 1400
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1450
$ws.Range("P3").Value = 1450
$ws.Range("D4").Value = 44477
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 1400
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 1460
$ws.Range("P4").Value = 1460
$ws.Range("D5").Value = 44868
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1300
$ws.Range("M5").Value = 1250
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 1250
$ws.Range("D6").Value = 44868
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 1000
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 1000
$ws.Range("D7").Value = 44468
$ws.Range("H7").Value = "Verde"
$ws.Range("K7").Value = 1800
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 1920
$ws.Range("P7").Value = 1920
$ws.Range("D9").Value = 44875
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 1500
$ws.Range("L9").Value = 1600
$ws.Range("M9").Value = 1550
$ws.Range("P9").Value = 1550
$ws.Range("D10").Value = 44519
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 1200
$ws.Range("L10").Value = 1300
$ws.Range("M10").Value = 1240
$ws.Range("P10").Value = 1240
$ws.Range("D11").Value = 44510
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 1300
$ws.Range("L11").Value = 1400
$ws.Range("M11").Value = 1350
$ws.Range("P11").Value = 1350
$ws.Range("D12").Value = 44860
$ws.Range("J12").Value = 1100
$ws.Range("K12").Value = 1500
$ws.Range("L12").Value = 1700
$ws.Range("M12").Value = 1609
$ws.Range("P12").Value = 1609
$ws.Range("D13").Value = 44876
$ws.Range("J13").Value = 350
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 1600
$ws.Range("M13").Value = 1557
$ws.Range("P13").Value = 1557
$ws.Range("D14").Value = 44511
$ws.Range("J14").Value = 600
$ws.Range("K14").Value = 1300
$ws.Range("L14").Value = 1400
$ws.Range("M14").Value = 1350
$ws.Range("P14").Value = 1350
$ws.Range("D15").Value = 44524
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 1500
$ws.Range("L15").Value = 1600
$ws.Range("M15").Value = 1550
$ws.Range("O15").Value = "Provincia de Talca"
$ws.Range("P15").Value = 1550
$ws.Range("D16").Value = 44881
$ws.Range("I16").Value = "Primera"
$ws.Range("K16").Value = 2600
$ws.Range("L16").Value = 2700
$ws.Range("M16").Value = 2650
$ws.Range("O16").Value = "Provincia de Linares"
$ws.Range("P16").Value = 2650
$ws.Range("D17").Value = 44881
$ws.Range("I17").Value = "Segunda"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 2400
$ws.Range("M17").Value = 2400
$ws.Range("P17").Value = 2400
$ws.Range("D19").Value = 44839
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("K19").Value = 1700
$ws.Range("L19").Value = 1800
$ws.Range("M19").Value = 1760
$ws.Range("P19").Value = 1760
